$wb = $excel.ActiveWorkbook

# --- Inventario: add "Sala" column (C) with the test room name ---
$wsInv = $wb.Worksheets.Item("Inventario")
$wsInv.Range("A1:A2").Copy()
$wsInv.Range("C1:C2").PasteSpecial(-4122)
$wsInv.Range("C1").Value = "Sala"
$wsInv.Range("C2").Value = "sala-de-prueba"

# --- Combinaciones: add "Sala" column (E) with the test room name ---
$wsComb = $wb.Worksheets.Item("Combinaciones")
$wsComb.Range("D1:D2").Copy()
$wsComb.Range("E1:E2").PasteSpecial(-4122)
$wsComb.Range("E1").Value = "Sala"
$wsComb.Range("E2").Value = "sala-de-prueba"

# --- Switch the active/selected tab from "sala-de-prueba" to "Combinaciones" ---
$wsComb.Activate()

$wb.Save()
